$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# IPL 2025 match results have come in - fill in the Match IDs for the
# matches that have now been played (rows 17-22, matches 16-21).

# Match 16 (row 17): Match ID rendered in Consolas, centered horizontally
# and vertically.
$d17 = $ws.Range("D17")
$d17.Value = 13388247143350
$d17.NumberFormat = "0"
$d17.Font.Name = "Consolas"
$d17.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$d17.VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter

# Match 17 (row 18)
$ws.Range("D18").Value = 13388319121419
$ws.Range("D18").NumberFormat = "0"

# Match 18 (row 19)
$ws.Range("D19").Value = 13388333622607
$ws.Range("D19").NumberFormat = "0"

# Match 19 (row 20)
$ws.Range("D20").Value = 13388419881809
$ws.Range("D20").NumberFormat = "0"

# Match 20 (row 21)
$ws.Range("D21").Value = 13388506308010
$ws.Range("D21").NumberFormat = "0"

# Match 21 (row 22)
$d22 = $ws.Range("D22")
$d22.Value = 13388578466393
$d22.NumberFormat = "0"

# Leave the selection where the author last clicked before saving
$d22.Select()
